$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("FoTOMRAEL")

# --- About sheet: add new "Mexico" annotation rows ---
$wsAbout.Range("A5").Value = "Mexico:"
$wsAbout.Range("A5").Font.Bold = $true

$wsAbout.Range("B5").Value = "Costs much more influenced by outside developments than the U.S. but wish to consider domestic policies. "

# B6 previously held an (empty) hyperlink-styled placeholder cell; clear that
# formatting before giving it its new plain-text note.
$wsAbout.Range("B6").Style = "Normal"
$wsAbout.Range("B6").Value = "Original US value 25% - changed to 50% for Mexico"

# The Hyperlink cell style is no longer used anywhere in the workbook, so drop it.
$wb.Styles.Item("Hyperlink").Delete()

# Minor column width tweaks
$wsAbout.Columns.Item(2).ColumnWidth = 51.6
$wsData.Columns.Item(1).ColumnWidth = 45.6

# --- FoTOMRAEL sheet: update the fraction for Mexico ---
$wsData.Range("B2").Value = 0.5
$wsData.Rows.Item(1).RowHeight = 16

# --- Restore cursor/selection state on each sheet ---
$wsData.Activate()
$wsData.Range("B3").Select()
$wsAbout.Activate()
$wsAbout.Range("A6").Select()
